$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column N is a new "2022" data column, mirroring the existing M (2021) column's
# formatting for each row. We copy the format from the corresponding M-column
# cell and then apply the value (and, where the source cell lacked the "0.0"
# number format, apply it) so the new cells line up with the rest of the table.

# Row 2 - blank separator row, just needs the thick-bottom-border formatting.
$ws.Range("M2").Copy() | Out-Null
$ws.Range("N2").PasteSpecial(-4122) | Out-Null

# Row 3 - year header.
$ws.Range("M3").Copy() | Out-Null
$ws.Range("N3").PasteSpecial(-4122) | Out-Null
$ws.Range("N3").Value = 2022

# Row 4 - bold summary row (already numeric-formatted like M4).
$ws.Range("M4").Copy() | Out-Null
$ws.Range("N4").PasteSpecial(-4122) | Out-Null
$ws.Range("N4").Value = 9.224468514531754

# Row 5
$ws.Range("M5").Copy() | Out-Null
$ws.Range("N5").PasteSpecial(-4122) | Out-Null
$ws.Range("N5").NumberFormat = "0.0"
$ws.Range("N5").Value = 4.6068543125097872

# Row 6
$ws.Range("M6").Copy() | Out-Null
$ws.Range("N6").PasteSpecial(-4122) | Out-Null
$ws.Range("N6").NumberFormat = "0.0"
$ws.Range("N6").Value = 13.543910285971602

# Row 7 - bold summary row.
$ws.Range("M7").Copy() | Out-Null
$ws.Range("N7").PasteSpecial(-4122) | Out-Null
$ws.Range("N7").NumberFormat = "0.0"
$ws.Range("N7").Value = 24.703327617190443

# Row 8
$ws.Range("M8").Copy() | Out-Null
$ws.Range("N8").PasteSpecial(-4122) | Out-Null
$ws.Range("N8").NumberFormat = "0.0"
$ws.Range("N8").Value = 28.608474183838851

# Row 9
$ws.Range("M9").Copy() | Out-Null
$ws.Range("N9").PasteSpecial(-4122) | Out-Null
$ws.Range("N9").NumberFormat = "0.0"
$ws.Range("N9").Value = 20.904451081350146

# Row 10 - bold summary row.
$ws.Range("M10").Copy() | Out-Null
$ws.Range("N10").PasteSpecial(-4122) | Out-Null
$ws.Range("N10").NumberFormat = "0.0"
$ws.Range("N10").Value = 26.720095429750884

# Row 11
$ws.Range("M11").Copy() | Out-Null
$ws.Range("N11").PasteSpecial(-4122) | Out-Null
$ws.Range("N11").NumberFormat = "0.0"
$ws.Range("N11").Value = 27.704327204727914

# Row 12 - bottom row with the thick bottom border.
$ws.Range("M12").Copy() | Out-Null
$ws.Range("N12").PasteSpecial(-4122) | Out-Null
$ws.Range("N12").NumberFormat = "0.0"
$ws.Range("N12").Value = 25.731792255708452

$excel.CutCopyMode = $false

$ws.Range("Q5").Select() | Out-Null
